# ============================================================================
# CompStat weekly report refresh: new reporting week + updated crime tallies
# ============================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: bump the report volume/issue number and the covered date range ---
$ws.Range("A8").Value = "Volume 30   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/30/2023  Through  2/5/2023"

# --- Stable donor cells used purely to restore the canonical cell style (cellXfs)
#     index after a value-type change. Excel/iron_native mint a brand new style
#     whenever a cell flips between number <-> text, so we copy-format back from a
#     cell elsewhere in the sheet that already carries the style we want to land on.
$styleDonorText = "A14"  # style used for "0"/"***.*" text placeholder cells
$styleDonorInt  = "C36"  # style used for whole-number cells (#,##0)
$styleDonorDec  = "K36"  # style used for one-decimal-place cells (#,##0.0)

function Set-TextPlaceholder($ref, $text) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $text
    $ws.Range($styleDonorText).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

function Set-NumberFromText($ref, $donor, $number) {
    $ws.Range($donor).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
    $ws.Range($ref).Value = $number
}

# ---- Row 15 ----
Set-TextPlaceholder "F15" "0"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = -100
Set-NumberFromText "M15" $styleDonorDec 0

# ---- Row 16 ----
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -42.105263157894
$ws.Range("I16").Value = 17
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = -15
$ws.Range("L16").Value = -10.526315789473
$ws.Range("M16").Value = -19.047619047619
$ws.Range("N16").Value = -81.521739130434

# ---- Row 17 ----
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 12
$ws.Range("J17").Value = 8
$ws.Range("L17").Value = 71.428571428571
$ws.Range("M17").Value = 71.428571428571

# ---- Row 18 ----
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = 38.888888888888
$ws.Range("I18").Value = 35
$ws.Range("J18").Value = 27
$ws.Range("K18").Value = 29.629629629629
$ws.Range("L18").Value = 66.666666666666
$ws.Range("M18").Value = 25
$ws.Range("N18").Value = -51.388888888888

# ---- Row 19 ----
$ws.Range("C19").Value = 26
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 73.333333333333
$ws.Range("F19").Value = 105
$ws.Range("G19").Value = 69
$ws.Range("H19").Value = 52.173913043478
$ws.Range("I19").Value = 137
$ws.Range("J19").Value = 94
$ws.Range("K19").Value = 45.744680851063
$ws.Range("L19").Value = 197.826086956522
$ws.Range("M19").Value = 52.222222222222
$ws.Range("N19").Value = -41.452991452991

# ---- Row 20 ----
Set-TextPlaceholder "C20" "0"
$ws.Range("F20").Value = 2
$ws.Range("H20").Value = -66.666666666666
$ws.Range("M20").Value = -25
$ws.Range("N20").Value = -96.511627906976

# ---- Row 21 ----
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 63.636363636363
$ws.Range("F21").Value = 148
$ws.Range("G21").Value = 118
$ws.Range("H21").Value = 25.423728813559
$ws.Range("I21").Value = 205
$ws.Range("J21").Value = 157
$ws.Range("K21").Value = 30.573248407643
$ws.Range("L21").Value = 118.085106382979
$ws.Range("M21").Value = 35.761589403973
$ws.Range("N21").Value = -59.405940594059

# ---- Row 22 ----
$ws.Range("C22").Value = 2
Set-TextPlaceholder "D22" "0"
Set-TextPlaceholder "E22" "***.*"
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 25
$ws.Range("I22").Value = 5
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 66.666666666666
$ws.Range("M22").Value = -58.333333333333

# ---- Row 24 ----
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 55
$ws.Range("F24").Value = 148
$ws.Range("G24").Value = 130
$ws.Range("H24").Value = 13.846153846153
$ws.Range("I24").Value = 186
$ws.Range("J24").Value = 160
$ws.Range("K24").Value = 16.25
$ws.Range("L24").Value = 60.344827586206
$ws.Range("M24").Value = 28.275862068965

# ---- Row 25 ----
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 60
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 92.857142857142
$ws.Range("I25").Value = 39
$ws.Range("J25").Value = 17
$ws.Range("K25").Value = 129.411764705882
$ws.Range("L25").Value = 105.263157894737
$ws.Range("M25").Value = 160

# ---- Row 26 ----
Set-TextPlaceholder "F26" "0"
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = -100

# ---- Row 27 ----
$ws.Range("C27").Value = 1
Set-NumberFromText "D27" $styleDonorInt 1
Set-NumberFromText "E27" $styleDonorDec 0
$ws.Range("I27").Value = 9
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 200
$ws.Range("L27").Value = 125

$excel.CutCopyMode = $false